# Updates cryptos list values (Price / Volume(1h)) per the commit diff.
# Values are written with a leading apostrophe to force text storage
# (Excel COM auto-converts plain numeric-looking strings like "604.89"
# into numbers, which would change the cell type away from the original
# inline-string/text representation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.667.95"
$ws.Range("E2").Value = "'  +1.19%  "

$ws.Range("D3").Value = "'3.303.26"
$ws.Range("E3").Value = "'  +5.90%  "

$ws.Range("E4").Value = "'  -0.10%  "

$ws.Range("D5").Value = "'604.89"
$ws.Range("E5").Value = "'  +1.90%  "

$ws.Range("D6").Value = "'142.25"
$ws.Range("E6").Value = "'  +4.72%  "

$ws.Range("E7").Value = "'  -0.08%  "

$ws.Range("D8").Value = "'3.298.76"
$ws.Range("E8").Value = "'  +5.91%  "

$ws.Range("E9").Value = "'  +1.05%  "

$ws.Range("E10").Value = "'  +3.52%  "

$ws.Range("D11").Value = "'5.49"
$ws.Range("E11").Value = "'  +4.89%  "

$ws.Range("E12").Value = "'  +3.95%  "

$ws.Range("E13").Value = "'  +1.54%  "

$ws.Range("D14").Value = "'34.67"
$ws.Range("E14").Value = "'  +2.09%  "

$ws.Range("D15").Value = "'3.841.54"

$ws.Range("E16").Value = "'  +1.26%  "

$ws.Range("D17").Value = "'3.298.61"
$ws.Range("E17").Value = "'  +6.01%  "

$ws.Range("D18").Value = "'63.753.45"
$ws.Range("E18").Value = "'  +1.18%  "

$ws.Range("E19").Value = "'  +3.71%  "

$ws.Range("D20").Value = "'480.72"
$ws.Range("E20").Value = "'  +2.39%  "

$ws.Range("D21").Value = "'14.13"
$ws.Range("E21").Value = "'  +0.61%  "

$ws.Range("E22").Value = "'  +5.66%  "

$ws.Range("D23").Value = "'8.03"
$ws.Range("E23").Value = "'  +5.34%  "

$ws.Range("D24").Value = "'13.50"
$ws.Range("E24").Value = "'  +4.81%  "

$ws.Range("D25").Value = "'84.75"
$ws.Range("E25").Value = "'  -0.84%  "

$ws.Range("E26").Value = "'  +0.08%  "

$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "'  +2.18%  "

$ws.Range("D28").Value = "'7.30"

$ws.Range("E29").Value = "'  -0.14%  "

$ws.Range("D30").Value = "'8.14"
$ws.Range("E30").Value = "'  +4.22%  "

$ws.Range("D31").Value = "'2.17"
$ws.Range("E31").Value = "'  +4.34%  "

$ws.Range("D32").Value = "'29.23"
$ws.Range("E32").Value = "'  +10.17%  "

$ws.Range("E33").Value = "'  -1.90%  "

$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "'  +1.01%  "

$ws.Range("D35").Value = "'1.09"
$ws.Range("E35").Value = "'  +2.54%  "

$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "'  +4.18%  "

$ws.Range("D37").Value = "'52.89"
$ws.Range("E37").Value = "'  +1.91%  "

$ws.Range("D38").Value = "'0.0₃0747"
$ws.Range("E38").Value = "'  +8.42%  "

$ws.Range("D39").Value = "'0.0400"
$ws.Range("E39").Value = "'  +3.93%  "

$ws.Range("D40").Value = "'425.97"
$ws.Range("E40").Value = "'  +2.56%  "

$ws.Range("D41").Value = "'3.055.09"
$ws.Range("E41").Value = "'  +5.56%  "

$ws.Range("E42").Value = "'  +2.55%  "

$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "'  +3.68%  "

$ws.Range("E44").Value = "'  -0.93%  "

$ws.Range("E45").Value = "'  +2.21%  "

$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "'  +5.39%  "

$ws.Range("D47").Value = "'26.31"
$ws.Range("E47").Value = "'  +4.14%  "

$ws.Range("E49").Value = "'  +2.49%  "

$ws.Range("E50").Value = "'  +3.18%  "

$ws.Range("D51").Value = "'124.17"
$ws.Range("E51").Value = "'  +3.31%  "
